$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:A27")
$rng.Sort($ws.Range("A1:A27"), 1) | Out-Null
$ws.Range("A26").Select() | Out-Null
